# Rename 'Codelists' sheet to 'Cells' and update its active selection,
# per commit "Rename 'codelists' to 'cells' / Close #256".

$wb = $excel.ActiveWorkbook

# Sheet4 ("Codelists") -> "Cells"
$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Make it the active sheet and move the selection from I17 to G18,
# matching the <selection activeCell="G18" sqref="G18"/> in the diff.
$ws.Activate()
$ws.Range("G18").Select()
